$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (old D -> E)
$ws.Columns.Item(4).Insert()

# Update header row
$ws.Range("C1").Value = "select"
$ws.Range("D1").Value = "select"
$ws.Range("E1").Value = "wait"

# Update row 3 (new D content is typed before B2's new hyperlink text)
$ws.Range("C3").Value = "Fiat"
$ws.Range("D3").Value = '{"target":"id=company","value":"Fiat"}'
$ws.Range("E3").Value = 2000

# Update row 2
$ws.Range("B2").Value = "http://127.0.0.1:9001/select/"
$ws.Range("C2").Value = "id=company"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""

# Set column widths to match Excel's "best fit" recalculation after edits
# (runtime stores ColumnWidth quantized to 1/7 + 5/7 padding, so pre-subtract padding)
$ws.Columns.Item(2).ColumnWidth = 27.875 - 0.7142857142857143
$ws.Columns.Item(3).ColumnWidth = 12.5 - 0.7142857142857143
$ws.Columns.Item(4).ColumnWidth = 36.5 - 0.7142857142857143
$ws.Columns.Item(5).ColumnWidth = 6.25 - 0.7142857142857143

# Selection
$ws.Range("D6").Select()
